$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New default row: eng / MOR / globaladmin / TRUE
$ws.Range("A2").Value = "eng"
$ws.Range("B2").Value = "MOR"
$ws.Range("C2").Value = "globaladmin"
$ws.Range("D2").Value = $true
$ws.Range("D2").NumberFormat = '"TRUE";"TRUE";"FALSE"'

# New default row: fra / MOR / globaladmin / TRUE
$ws.Range("A3").Value = "fra"
$ws.Range("B3").Value = "MOR"
$ws.Range("C3").Value = "globaladmin"
$ws.Range("D3").Value = $true
$ws.Range("D3").NumberFormat = '"TRUE";"TRUE";"FALSE"'

# Column widths for zone_code / usr_id columns
$ws.Columns.Item(2).ColumnWidth = 12
$ws.Columns.Item(3).ColumnWidth = 9.8

# Restore selection to match the committed state
$null = $ws.Range("B10").Select()
